$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 33-36: the dependency columns (E/F) shift down
# by one "type" (a new Server node is threaded in above them at row 32).
# Row 32 itself is fixed up last, below, once its new dependency strings
# already exist (matches how the new shared strings line up in the
# original edit).
$ws.Cells.Item(33, 5).Value = "Facilities"
$ws.Cells.Item(33, 6).Value = "Location 20"

$ws.Cells.Item(34, 5).Value = "Procurements"
$ws.Cells.Item(34, 6).Value = "PO 20"

$ws.Cells.Item(35, 5).Value = "People"
$ws.Cells.Item(35, 6).Value = "Person 20"

$ws.Cells.Item(36, 5).Value = "Data"
$ws.Cells.Item(36, 6).Value = "Data 20"

# --- Insert one new row at 42 so the trailing blank formatting row moves
# from row 42 down to row 43, opening up rows 37-41 for the new data.
$ws.Rows.Item(42).Insert(-4121)

# --- Clone row 36 (values + formats) into each of the newly opened rows
# 37-41 individually. Columns C/G keep the long "Business description..."
# text and the wrapped style (s=4) this way; doing it via two passes (a
# full paste, then a formats-only re-paste) avoids this host's auto-height
# recompute that an outright .Value write into a wrap-styled cell would
# otherwise trigger, so rows 37-41 keep the plain 71.25pt row height
# instead of a custom one. (Pasting row-by-row, rather than into the whole
# A37:G41 block at once, keeps every column's value intact on every row.)
$ws.Range("A36:G36").Copy()
$ws.Range("A37:G37").PasteSpecial(-4104)
$ws.Range("A38:G38").PasteSpecial(-4104)
$ws.Range("A39:G39").PasteSpecial(-4104)
$ws.Range("A40:G40").PasteSpecial(-4104)
$ws.Range("A41:G41").PasteSpecial(-4104)

$ws.Range("A36:G36").Copy()
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A38:G38").PasteSpecial(-4122)
$ws.Range("A39:G39").PasteSpecial(-4122)
$ws.Range("A40:G40").PasteSpecial(-4122)
$ws.Range("A41:G41").PasteSpecial(-4122)

# Row 37: Applications / App 2 -> Technology / Tech 20
$ws.Cells.Item(37, 5).Value = "Technology"
$ws.Cells.Item(37, 6).Value = "Tech 20"

# Row 38: People / People 1 -> Procurements / PO 1
$ws.Cells.Item(38, 1).Value = "People"
$ws.Cells.Item(38, 2).Value = "People 1"
$ws.Cells.Item(38, 5).Value = "Procurements"
$ws.Cells.Item(38, 6).Value = "PO 1"

# Row 39: People / People 1 -> Applications / App 3
$ws.Cells.Item(39, 1).Value = "People"
$ws.Cells.Item(39, 2).Value = "People 1"
$ws.Cells.Item(39, 5).Value = "Applications"
$ws.Cells.Item(39, 6).Value = "App 3"

# Row 40: People / People 1 -> Facilities / Location 1
$ws.Cells.Item(40, 1).Value = "People"
$ws.Cells.Item(40, 2).Value = "People 1"
$ws.Cells.Item(40, 5).Value = "Facilities"
$ws.Cells.Item(40, 6).Value = "Location 1"

# Row 41: Facilities / Location 1 -> Server / Server 1 (the new Server
# node itself, depending on the facility it's racked in).
$ws.Cells.Item(41, 1).Value = "Facilities"
$ws.Cells.Item(41, 2).Value = "Location 1"
$ws.Cells.Item(41, 5).Value = "Server"
$ws.Cells.Item(41, 6).Value = "Server 1"

# --- Finally, row 32: CI_Name moves from "App 2" to "OIT" with its
# dependency becoming the new Server node ("Server 1") just added above.
$ws.Cells.Item(32, 2).Value = "OIT"
$ws.Cells.Item(32, 5).Value = "Server "
$ws.Cells.Item(32, 6).Value = "Server 1"

# --- Sheet view: selection on L33 (scrolled down to show the new rows in
# real Excel via topLeftCell, which this headless host does not expose).
$ws.Range("L33").Select()
